$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "PickAndPlace_PCB V1.1_No_ANT_20"

# Update C3 row (row 9): Mid X/Y, Ref X/Y, Pad X/Y
$ws.Range("D9").Value = "60.198mm"
$ws.Range("E9").Value = "45.212mm"
$ws.Range("F9").Value = "60.198mm"
$ws.Range("G9").Value = "45.212mm"
$ws.Range("H9").Value = "60.198mm"
$ws.Range("I9").Value = "46.212mm"

# Update C4 row (row 16): Mid X/Y, Ref X/Y, Pad X/Y
$ws.Range("D16").Value = "62.357mm"
$ws.Range("E16").Value = "45.212mm"
$ws.Range("F16").Value = "62.357mm"
$ws.Range("G16").Value = "45.212mm"
$ws.Range("H16").Value = "62.357mm"
$ws.Range("I16").Value = "46.212mm"

# Add new row 18 for J3 (SMD Jumper 3-pin)
$ws.Range("A18").Value = "J3"
$ws.Range("B18").Value = "SMD_JUMPER_3PIN"
$ws.Range("C18").Value = "SDM_JUMPER_3PIN"
$ws.Range("D18").Value = "61.686mm"
$ws.Range("E18").Value = "41.484mm"
$ws.Range("F18").Value = "61.686mm"
$ws.Range("G18").Value = "41.484mm"
$ws.Range("H18").Value = "62.686mm"
$ws.Range("I18").Value = "41.484mm"
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = "T"
$ws.Range("L18").Value = 180
$ws.Range("M18").Value = "Yes"
$ws.Range("N18").Value = "SMD_JUMPER_3PIN"
